$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-7 (TruckID, AssignedDockPosition, start_loading_time, end_loading_time)
$data = @(
    @(4, 3, 5, 5),
    @(1, 4, 5, 6),
    @(5, 5, 5, 6),
    @(5, 5, 11, 12),
    @(2, 7, 5, 6),
    @(3, 7, 11, 12)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
}
